# Insert a new row at position 386, shifting existing rows 386-471 down to 387-472.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(386).Insert()

# Populate the newly inserted row 386 with the new weekly data point.
$ws.Cells.Item(386, 1).Value2  = 6
$ws.Cells.Item(386, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(386, 3).Value2  = "Metropolitana"
$ws.Cells.Item(386, 4).Value2  = 44641
$ws.Cells.Item(386, 5).Value2  = 13
$ws.Cells.Item(386, 6).Value2  = 100112044
$ws.Cells.Item(386, 7).Value2  = "Perejil"
$ws.Cells.Item(386, 8).Value2  = "Sin especificar"
$ws.Cells.Item(386, 9).Value2  = "Primera"
$ws.Cells.Item(386, 10).Value2 = 150
$ws.Cells.Item(386, 11).Value2 = 15000
$ws.Cells.Item(386, 12).Value2 = 16000
$ws.Cells.Item(386, 13).Value2 = 15400
$ws.Cells.Item(386, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(386, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(386, 16).Value2 = 5133
$ws.Cells.Item(386, 17).Value2 = 3
$ws.Cells.Item(386, 18).Value2 = "Hortaliza"
